$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B:E are treated as text so values like "1.00" or "61.563.22" are
# not reinterpreted as numbers/dates by Excel, matching the source inline-string data.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "61.563.22"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "3.449.50"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "578.88"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").Value = "144.44"
$ws.Range("E6").Value = "  +6.51%  "
$ws.Range("D7").Value = "3.451.27"
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.476"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("E11").Value = "  +3.50%  "
$ws.Range("E12").Value = "  +2.26%  "
$ws.Range("D13").Value = "4.042.70"
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("D14").Value = "28.04"
$ws.Range("E14").Value = "  +9.89%  "
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.453.08"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000173"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").Value = "61.735.49"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "6.26"
$ws.Range("E19").Value = "  +8.51%  "
$ws.Range("D20").Value = "14.23"
$ws.Range("E20").Value = "  +4.00%  "
$ws.Range("D21").Value = "9.51"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").Value = "389.02"
$ws.Range("E22").Value = "  +4.57%  "
$ws.Range("D23").Value = "0.563"
$ws.Range("E23").Value = "  +3.52%  "
$ws.Range("D24").Value = "73.40"
$ws.Range("E24").Value = "  +3.80%  "
$ws.Range("D25").Value = "5.78"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "0.0000123"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").Value = "3.587.93"
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("D29").Value = "0.179"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").Value = "7.64"
$ws.Range("E30").Value = "  +4.41%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.47"
$ws.Range("E32").Value = "  -9.84%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "8.14"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "2.18"
$ws.Range("E34").Value = "  +2.88%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "24.02"
$ws.Range("E36").Value = "  +3.68%  "
$ws.Range("D37").Value = "3.480.74"
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("D38").Value = "6.99"
$ws.Range("E38").Value = "  +3.36%  "
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").Value = "166.99"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("D42").Value = "28.12"
$ws.Range("E42").Value = "  +12.91%  "
$ws.Range("D43").Value = "0.0780"
$ws.Range("E43").Value = "  +3.74%  "
$ws.Range("D44").Value = "0.802"
$ws.Range("E44").Value = "  +4.09%  "
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.73"
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "42.29"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("D48").Value = "4.47"
$ws.Range("E48").Value = "  +4.22%  "
$ws.Range("D49").Value = "2.584.59"
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("E51").Value = "  +2.49%  "
